$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 202.93333
$ws.Range("I28").Value = 205.07692
$ws.Range("J28").Value = 189
$ws.Range("K28").Value = 205.07692
$ws.Range("L28").Value = 189
$ws.Range("M28").Value = 279.92308
$ws.Range("N28").Value = -1159
$ws.Range("H115").Value = 1177
$ws.Range("I115").Value = 721.25
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 2163.75
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -596.75
$ws.Range("N115").Value = -12134
$ws.Range("H116").Value = 698465.7
$ws.Range("J116").Value = 5122.1816
$ws.Range("L116").Value = 5122.1816
$ws.Range("N116").Value = -12006.1816
$ws.Range("H118").Value = 903.4545000000001
$ws.Range("I118").Value = 619.75
$ws.Range("J118").Value = 1065.5714
$ws.Range("K118").Value = 1859.25
$ws.Range("L118").Value = 3196.7142
$ws.Range("M118").Value = -202.25
$ws.Range("N118").Value = -6510.7142
$ws.Range("H132").Value = 3032.1428
$ws.Range("I132").Value = 2874.2654
$ws.Range("K132").Value = 8622.796200000001
$ws.Range("M132").Value = -6092.796200000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8727.687
$ws.Range("I32").Value = 8675.754999999999
$ws.Range("K32").Value = 8675.754999999999
$ws.Range("M32").Value = -8388.754999999999
$ws.Range("H61").Value = 9476
$ws.Range("I61").Value = 11538.723
$ws.Range("J61").Value = 6100.636
$ws.Range("K61").Value = 11538.723
$ws.Range("L61").Value = 6100.636
$ws.Range("M61").Value = -11326.723
$ws.Range("N61").Value = -6524.636
$ws.Range("H136").Value = 9476
$ws.Range("I136").Value = 11538.723
$ws.Range("J136").Value = 6100.636
$ws.Range("K136").Value = 34616.169
$ws.Range("L136").Value = 18301.908
$ws.Range("M136").Value = -32066.169
$ws.Range("N136").Value = -23401.908

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3122.1304
$ws.Range("I20").Value = 1723.3572
$ws.Range("K20").Value = 1723.3572
$ws.Range("M20").Value = -1476.3572
$ws.Range("H80").Value = 228.84
$ws.Range("J80").Value = 254.05263
$ws.Range("L80").Value = 254.05263
$ws.Range("N80").Value = -2250.05263
$ws.Range("H83").Value = 228.84
$ws.Range("J83").Value = 254.05263
$ws.Range("L83").Value = 1270.26315
$ws.Range("N83").Value = -11254.26315
$ws.Range("H134").Value = 10486.421
$ws.Range("J134").Value = 6398
$ws.Range("L134").Value = 19194
$ws.Range("N134").Value = -24264

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 200.5
$ws.Range("I7").Value = 192.33333
$ws.Range("J7").Value = 225
$ws.Range("K7").Value = 192.33333
$ws.Range("L7").Value = 225
$ws.Range("M7").Value = -79.33332999999999
$ws.Range("N7").Value = -451
$ws.Range("H16").Value = 1310
$ws.Range("J16").Value = 1310
$ws.Range("L16").Value = 1310
$ws.Range("N16").Value = -1884
$ws.Range("H31").Value = 4813.227
$ws.Range("I31").Value = 4879.543
$ws.Range("K31").Value = 4879.543
$ws.Range("M31").Value = -4584.543
$ws.Range("H34").Value = 4813.227
$ws.Range("I34").Value = 4879.543
$ws.Range("K34").Value = 4879.543
$ws.Range("M34").Value = -4677.543
$ws.Range("H41").Value = 312.5
$ws.Range("I41").Value = 312.5
$ws.Range("K41").Value = 312.5
$ws.Range("M41").Value = 115.5
$ws.Range("H99").Value = 316888.38
$ws.Range("I99").Value = 558023.75
$ws.Range("K99").Value = 558023.75
$ws.Range("M99").Value = -556525.75
$ws.Range("H105").Value = 8349
$ws.Range("I105").Value = 11260.8
$ws.Range("K105").Value = 11260.8
$ws.Range("M105").Value = -9513.799999999999
$ws.Range("H113").Value = 1310
$ws.Range("J113").Value = 1310
$ws.Range("L113").Value = 1310
$ws.Range("N113").Value = -5650
$ws.Range("H126").Value = 316888.38
$ws.Range("I126").Value = 558023.75
$ws.Range("K126").Value = 1674071.25
$ws.Range("M126").Value = -1671601.25
$ws.Range("H134").Value = 11610.667
$ws.Range("I134").Value = 14925.333
$ws.Range("K134").Value = 44775.999
$ws.Range("M134").Value = -42240.999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 299.8889
$ws.Range("I29").Value = 261.9
$ws.Range("J29").Value = 408.42856
$ws.Range("K29").Value = 785.6999999999999
$ws.Range("L29").Value = 1225.28568
$ws.Range("M29").Value = -508.6999999999999
$ws.Range("N29").Value = -1779.28568
$ws.Range("H116").Value = 6875
$ws.Range("I116").Value = 3750
$ws.Range("K116").Value = 11250
$ws.Range("M116").Value = -7808

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 430.92856
$ws.Range("J107").Value = 292
$ws.Range("L107").Value = 292
$ws.Range("N107").Value = -4132
$ws.Range("H113").Value = 11522.667
$ws.Range("I113").Value = 19212
$ws.Range("J113").Value = 3833.3333
$ws.Range("K113").Value = 19212
$ws.Range("L113").Value = 3833.3333
$ws.Range("M113").Value = -17042
$ws.Range("N113").Value = -8173.3333

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 20000
$ws.Range("I56").Value = 20000
$ws.Range("K56").Value = 20000
$ws.Range("M56").Value = -19309
$ws.Range("H61").Value = 4891.625
$ws.Range("I61").Value = 1186.2
$ws.Range("J61").Value = 11067.333
$ws.Range("K61").Value = 1186.2
$ws.Range("L61").Value = 11067.333
$ws.Range("M61").Value = -984.2
$ws.Range("N61").Value = -11471.333
$ws.Range("H82").Value = 2963.3845
$ws.Range("I82").Value = 4078.4285
$ws.Range("J82").Value = 1662.5
$ws.Range("K82").Value = 4078.4285
$ws.Range("L82").Value = 1662.5
$ws.Range("M82").Value = -3717.4285
$ws.Range("N82").Value = -2384.5
$ws.Range("H85").Value = 2963.3845
$ws.Range("I85").Value = 4078.4285
$ws.Range("J85").Value = 1662.5
$ws.Range("K85").Value = 4078.4285
$ws.Range("L85").Value = 1662.5
$ws.Range("M85").Value = -2830.4285
$ws.Range("N85").Value = -4158.5
$ws.Range("H113").Value = 4891.625
$ws.Range("I113").Value = 1186.2
$ws.Range("J113").Value = 11067.333
$ws.Range("K113").Value = 1186.2
$ws.Range("L113").Value = 11067.333
$ws.Range("M113").Value = 983.8
$ws.Range("N113").Value = -15407.333
$ws.Range("H122").Value = 6371.2173
$ws.Range("I122").Value = 5835.8667
$ws.Range("J122").Value = 7375
$ws.Range("K122").Value = 17507.6001
$ws.Range("L122").Value = 22125
$ws.Range("M122").Value = -15057.6001
$ws.Range("N122").Value = -27025
$ws.Range("H136").Value = 6214.125
$ws.Range("I136").Value = 3642.6
$ws.Range("J136").Value = 7383
$ws.Range("K136").Value = 10927.8
$ws.Range("L136").Value = 22149
$ws.Range("M136").Value = -8377.799999999999
$ws.Range("N136").Value = -27249

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 35330.11
$ws.Range("I107").Value = 2246.375
$ws.Range("J107").Value = 300000
$ws.Range("K107").Value = 6739.125
$ws.Range("L107").Value = 900000
$ws.Range("M107").Value = -4819.125
$ws.Range("N107").Value = -903840
$ws.Range("H132").Value = 17563.44
$ws.Range("I132").Value = 30917.363
$ws.Range("J132").Value = 7071.0713
$ws.Range("K132").Value = 92752.08900000001
$ws.Range("L132").Value = 21213.2139
$ws.Range("M132").Value = -90222.08900000001
$ws.Range("N132").Value = -26273.2139
